$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (from 12.5546875 to 13.5546875).
# Note: the host engine quantizes ColumnWidth to 1/6-character steps when
# converting through its internal points/pixel representation, so the exact
# fractional value cannot be reproduced bit-for-bit; 12.6667 lands on the
# closest reachable stored width (13.5) to the target (13.5546875).
$ws.Columns.Item(1).ColumnWidth = 12.6667

# Update values in column A (rows 1-59)
$values = @(
    0.448282122327467,
    0.70074058298358577,
    0.83881179074299539,
    0.76869366594304445,
    0.52353537537000716,
    0.71217793281009323,
    0.78937720966030911,
    0.87538088455738061,
    0.97247305901421544,
    0.8041620966481271,
    0.90978972120237123,
    0.49937253004245841,
    0.77146889666771634,
    0.73839037103639682,
    0.45701038074698336,
    0.92876594140953617,
    0.64579004771583426,
    0.622075438252644,
    0.70811853667873403,
    0.53424732674157149,
    0.33682978169698902,
    0.18632177870822045,
    0.2610117402044112,
    0.37745141298087515,
    0.35414441718317524,
    0.54934238734243079,
    0.77648500740391602,
    0.57953540007109883,
    0.42326334027283968,
    0.28345249948142059,
    0.22333621015638833,
    0.23091351448890512,
    0.31160704949718931,
    0.19175849765808459,
    0.39524850203792539,
    0.25150263430876751,
    0.26157829593852366,
    0.33866680274864958,
    0.33761284548147025,
    0.3188313926620947,
    0.15875592171552239,
    0.15092078513953033,
    0.071637956543694939,
    0.01488823162330823,
    0.04360888423032358,
    0.019856575224641411,
    0.09095679205209567,
    0.042329408498803407,
    0.078819285913480225,
    0.054454146531799544,
    0.092612123455437817,
    0.12098233248029841,
    0.18819693986188979,
    0.20884592664848661,
    0.193498032923139,
    0.29199399361980899,
    0.49427815822211307,
    0.25617376187097568,
    0.27982352138315891
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
